$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: worksheet row number, new Price (column D) text, new Volume(1h) text (column E).
# $null means "leave this cell unchanged" (that column did not change for this row).
# NumericD marks Price values that look like plain numbers ("7.97") so we must write them
# with a leading apostrophe to keep them stored as text instead of being parsed as numbers
# (matches values such as "66.666.35" which already fail numeric parsing and stay text).
$updates = @(
    @{ Row = 2; D = '66.666.35'; NumericD = $false; E = '  +0.53%  ' }
    @{ Row = 3; D = '3.570.71'; NumericD = $false; E = '  +0.48%  ' }
    @{ Row = 4; D = $null; NumericD = $false; E = '  -0.01%  ' }
    @{ Row = 5; D = '608.44'; NumericD = $true; E = '  +0.03%  ' }
    @{ Row = 6; D = '145.55'; NumericD = $true; E = '  +0.75%  ' }
    @{ Row = 7; D = $null; NumericD = $false; E = '  +0.15%  ' }
    @{ Row = 8; D = '0.492'; NumericD = $true; E = '  +2.35%  ' }
    @{ Row = 10; D = '7.97'; NumericD = $true; E = '  -1.42%  ' }
    @{ Row = 11; D = '0.417'; NumericD = $true; E = '  +1.31%  ' }
    @{ Row = 12; D = '4.174.19'; NumericD = $false; E = '  +0.46%  ' }
    @{ Row = 13; D = $null; NumericD = $false; E = '  +0.19%  ' }
    @{ Row = 14; D = '30.10'; NumericD = $true; E = '  -0.09%  ' }
    @{ Row = 15; D = '3.531.03'; NumericD = $false; E = '  -0.73%  ' }
    @{ Row = 16; D = '66.691.25'; NumericD = $false; E = '  +0.48%  ' }
    @{ Row = 17; D = $null; NumericD = $false; E = '  +0.13%  ' }
    @{ Row = 18; D = '11.39'; NumericD = $true; E = '  +2.32%  ' }
    @{ Row = 19; D = $null; NumericD = $false; E = '  +0.17%  ' }
    @{ Row = 20; D = '15.09'; NumericD = $true; E = '  +1.03%  ' }
    @{ Row = 21; D = '433.14'; NumericD = $true; E = '  +1.20%  ' }
    @{ Row = 22; D = '0.619'; NumericD = $true; E = '  +2.60%  ' }
    @{ Row = 23; D = $null; NumericD = $false; E = '  +0.45%  ' }
    @{ Row = 24; D = '3.708.41'; NumericD = $false; E = '  +0.44%  ' }
    @{ Row = 25; D = $null; NumericD = $false; E = '  +0.03%  ' }
    @{ Row = 26; D = $null; NumericD = $false; E = '  -1.46%  ' }
    @{ Row = 27; D = $null; NumericD = $false; E = '  -0.35%  ' }
    @{ Row = 28; D = $null; NumericD = $false; E = '  +1.01%  ' }
    @{ Row = 29; D = '9.21'; NumericD = $true; E = '  +0.20%  ' }
    @{ Row = 30; D = $null; NumericD = $false; E = '  -0.10%  ' }
    @{ Row = 31; D = '3.563.72'; NumericD = $false; E = '  +0.55%  ' }
    @{ Row = 32; D = $null; NumericD = $false; E = '  -2.56%  ' }
    @{ Row = 33; D = '25.43'; NumericD = $true; E = '  +0.11%  ' }
    @{ Row = 34; D = $null; NumericD = $false; E = '  -1.85%  ' }
    @{ Row = 35; D = '7.87'; NumericD = $true; E = '  +0.42%  ' }
    @{ Row = 36; D = $null; NumericD = $false; E = '  +0.03%  ' }
    @{ Row = 37; D = $null; NumericD = $false; E = '  -1.92%  ' }
    @{ Row = 38; D = $null; NumericD = $false; E = '  -0.12%  ' }
    @{ Row = 39; D = '172.82'; NumericD = $true; E = '  -1.07%  ' }
    @{ Row = 40; D = '0.0854'; NumericD = $true; E = '  -0.52%  ' }
    @{ Row = 41; D = '5.23'; NumericD = $true; E = '  -0.72%  ' }
    @{ Row = 42; D = $null; NumericD = $false; E = '  -0.45%  ' }
    @{ Row = 43; D = $null; NumericD = $false; E = '  +1.39%  ' }
    @{ Row = 44; D = $null; NumericD = $false; E = '  -0.02%  ' }
    @{ Row = 45; D = '2.53'; NumericD = $true; E = '  +5.44%  ' }
    @{ Row = 46; D = $null; NumericD = $false; E = '  -1.66%  ' }
    @{ Row = 47; D = '25.16'; NumericD = $true; E = '  -3.37%  ' }
    @{ Row = 48; D = '7.22'; NumericD = $true; E = $null }
    @{ Row = 49; D = '23.54'; NumericD = $true; E = '  +2.72%  ' }
    @{ Row = 50; D = '0.942'; NumericD = $true; E = $null }
    @{ Row = 51; D = $null; NumericD = $false; E = '  -1.15%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        if ($u.NumericD) {
            # Leading apostrophe forces Excel to keep the digits-and-dots text as a string
            $dCell.Formula = "'" + $u.D
        } else {
            $dCell.Value = $u.D
        }
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}
